$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 37; this shifts the existing rows 37..73 down to 38..74
$ws.Rows.Item(37).Insert()

# Populate the newly inserted row 37 with the new weekly record.
# Columns A,B,C,E,F,G,H,I,N,O,Q,R are constant across this sheet's data rows.
$ws.Cells.Item(37, 1).Value = 3
$ws.Cells.Item(37, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(37, 3).Value = "Coquimbo"
$ws.Cells.Item(37, 4).Value = 44790
$ws.Cells.Item(37, 5).Value = 5
$ws.Cells.Item(37, 6).Value = 100112035
$ws.Cells.Item(37, 7).Value = "Bruselas (repollito)"
$ws.Cells.Item(37, 8).Value = "Sin especificar"
$ws.Cells.Item(37, 9).Value = "Primera"
$ws.Cells.Item(37, 10).Value = 40
$ws.Cells.Item(37, 11).Value = 15000
$ws.Cells.Item(37, 12).Value = 15000
$ws.Cells.Item(37, 13).Value = 15000
$ws.Cells.Item(37, 14).Value = "$/malla 15 kilos"
$ws.Cells.Item(37, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(37, 16).Value = 1000
$ws.Cells.Item(37, 17).Value = 15
$ws.Cells.Item(37, 18).Value = "Hortaliza"
